$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Texas -- Bexar County ---
$ws.Range("C3").Value = 33555
$ws.Range("D3").Value = 283

# --- Row 4: New York -- New York ---
# Copy the date number format from an existing "Date Published" cell (B3)
# so B4 picks up the same style index used elsewhere in that column.
$ws.Range("B4").NumberFormat = $ws.Range("B3").NumberFormat
$ws.Range("B4").Value = 44034

# C4/D4 keep arriving as literal numeric-looking text (not true numbers) in
# this run's source data, so force Text via the apostrophe-prefix convention.
$ws.Range("C4").Value = "'219128"
$ws.Range("D4").Value = "'18803"

$ws.Range("E4").Value = 33790
$ws.Range("F4").Value = 5239
$ws.Range("G4").Value = 30.07
$ws.Range("H4").Value = 30.43

$ws.Range("J4").Value = $true

$ws.Range("K4").Value = 112360
$ws.Range("L4").Value = 17217

$ws.Range("O4").Value = "Success!"

# --- Row 11: California - San Diego ---
$ws.Range("B11").Value = 44034
$ws.Range("C11").Value = 25107
$ws.Range("D11").Value = 505
$ws.Range("E11").Value = 918
$ws.Range("F11").Value = 19
$ws.Range("G11").Value = 4.67
$ws.Range("H11").Value = 3.85
$ws.Range("K11").Value = 19655
$ws.Range("L11").Value = 493
